$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" title text (A1)
$ws.Range("A1").Value = "Datos actualizados a 25 de Mayo de 2020 a las 18:35"

# Swap the province labels for rows 13 and 14 (Aragon <-> Valencia/Valencia)
$ws.Range("A13").Value = "Valencia/Valencia"
$ws.Range("A14").Value = "Aragon"

# Update numeric data per row (Casos totales=B, Casos activos=C, Recuperados=D, Muertes=E)
# Row 4 - Madrid
$ws.Range("B4").Value = 67932
$ws.Range("D4").Value = 18510
$ws.Range("E4").Value = 8686

# Row 5 - Cataluña
$ws.Range("B5").Value = 57320
$ws.Range("D5").Value = 25542
$ws.Range("E5").Value = 5575

# Row 6 - Castilla y Leon
$ws.Range("B6").Value = 18555
$ws.Range("D6").Value = 7926
$ws.Range("E6").Value = 1913

# Row 7 - Castilla-La Mancha
$ws.Range("B7").Value = 16909
$ws.Range("D7").Value = 7729
$ws.Range("E7").Value = 2788

# Row 9 - Andalucia
$ws.Range("B9").Value = 12450
$ws.Range("D9").Value = 445
$ws.Range("E9").Value = 1334

# Row 13 - now Valencia/Valencia (values below come straight from the diff for row 13)
$ws.Range("B13").Value = 5609
$ws.Range("C13").Value = 4907
$ws.Range("D13").Value = 2767
$ws.Range("E13").Value = 693

# Row 14 - now Aragon (values below come straight from the diff for row 14)
$ws.Range("B14").Value = 5600
$ws.Range("C14").Value = 3772
$ws.Range("D14").Value = 1002
$ws.Range("E14").Value = 826

# Row 16 - Navarra
$ws.Range("B16").Value = 5195
$ws.Range("D16").Value = 800
$ws.Range("E16").Value = 490

# Row 20 - La Rioja
$ws.Range("B20").Value = 4041
$ws.Range("D20").Value = 576

# Row 32 - Asturias
$ws.Range("B32").Value = 2380
$ws.Range("D32").Value = 1029
$ws.Range("E32").Value = 288

# Row 33 - Gran Canaria
$ws.Range("B33").Value = 2317
$ws.Range("D33").Value = 630
$ws.Range("E33").Value = 150

# Row 40 - Murcia
$ws.Range("B40").Value = 1573
$ws.Range("E40").Value = 139
